# Auto-generated edit script: applies the Typhon_Profits value updates
# (the workbook stores each in-game sheet as a separate worksheet named
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR; the source diff concatenated all 8
# sheets XML one after another, so row numbers repeat per sheet).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 483.85715
$ws.Range("I38").Value = 99
$ws.Range("J38").Value = 772.5
$ws.Range("K38").Value = 297
$ws.Range("L38").Value = 2317.5
$ws.Range("M38").Value = 75
$ws.Range("N38").Value = -3061.5
$ws.Range("H51").Value = 5555.6665
$ws.Range("I51").Value = 6999.75
$ws.Range("K51").Value = 6999.75
$ws.Range("M51").Value = -6515.75
$ws.Range("H94").Value = 2750
$ws.Range("I94").Value = 2750
$ws.Range("K94").Value = 2750
$ws.Range("M94").Value = -2299
$ws.Range("H112").Value = 5209466
$ws.Range("J112").Value = 1181.6957
$ws.Range("L112").Value = 3545.0871
$ws.Range("N112").Value = -5761.0871
$ws.Range("H118").Value = 775
$ws.Range("J118").Value = 1200
$ws.Range("L118").Value = 3600
$ws.Range("N118").Value = -6914
$ws.Range("H135").Value = 21745628
$ws.Range("I135").Value = 1120.2667
$ws.Range("J135").Value = 62516580
$ws.Range("K135").Value = 10082.4003
$ws.Range("L135").Value = 562649220
$ws.Range("M135").Value = -7547.400299999999
$ws.Range("N135").Value = -562654290
$ws.Range("H137").Value = 114013.055
$ws.Range("I137").Value = 156120.34
$ws.Range("K137").Value = 468361.02
$ws.Range("M137").Value = -465811.02
$ws.Range("H138").Value = 3921.4648
$ws.Range("J138").Value = 3731.1846
$ws.Range("L138").Value = 11193.5538
$ws.Range("N138").Value = -21473.5538

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -827
$ws.Range("H32").Value = 13109.444
$ws.Range("I32").Value = 9612.846
$ws.Range("J32").Value = 29638.818
$ws.Range("K32").Value = 9612.846
$ws.Range("L32").Value = 29638.818
$ws.Range("M32").Value = -9325.846
$ws.Range("N32").Value = -30212.818
$ws.Range("H45").Value = 4648.421
$ws.Range("I45").Value = 5159.2
$ws.Range("K45").Value = 5159.2
$ws.Range("M45").Value = -4782.2
$ws.Range("H61").Value = 35139996
$ws.Range("I61").Value = 87843336
$ws.Range("J61").Value = 4433.3335
$ws.Range("K61").Value = 87843336
$ws.Range("L61").Value = 4433.3335
$ws.Range("M61").Value = -87843124
$ws.Range("N61").Value = -4857.3335
$ws.Range("I88").Value = 1850
$ws.Range("J88").Value = 168110.17
$ws.Range("K88").Value = 1850
$ws.Range("L88").Value = 168110.17
$ws.Range("M88").Value = -1444
$ws.Range("N88").Value = -168922.17
$ws.Range("I91").Value = 1850
$ws.Range("J91").Value = 168110.17
$ws.Range("K91").Value = 1850
$ws.Range("L91").Value = 168110.17
$ws.Range("M91").Value = -446
$ws.Range("N91").Value = -170918.17
$ws.Range("H132").Value = 18540836
$ws.Range("I132").Value = 31253564
$ws.Range("K132").Value = 93760692
$ws.Range("M132").Value = -93758162
$ws.Range("H136").Value = 35139996
$ws.Range("I136").Value = 87843336
$ws.Range("J136").Value = 4433.3335
$ws.Range("K136").Value = 263530008
$ws.Range("L136").Value = 13300.0005
$ws.Range("M136").Value = -263527458
$ws.Range("N136").Value = -18400.0005

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2062.5454
$ws.Range("I3").Value = 2012.5714
$ws.Range("K3").Value = 2012.5714
$ws.Range("M3").Value = -1898.5714
$ws.Range("H86").Value = 2056.2354
$ws.Range("I86").Value = 1869.9333
$ws.Range("K86").Value = 1869.9333
$ws.Range("M86").Value = -746.9332999999999
$ws.Range("H89").Value = 2056.2354
$ws.Range("I89").Value = 1869.9333
$ws.Range("K89").Value = 9349.666499999999
$ws.Range("M89").Value = -3733.666499999999
$ws.Range("H134").Value = 3523.18
$ws.Range("I134").Value = 3183.0466
$ws.Range("K134").Value = 9549.139800000001
$ws.Range("M134").Value = -7014.139800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6120.123
$ws.Range("I31").Value = 3496.875
$ws.Range("J31").Value = 7143.829
$ws.Range("K31").Value = 3496.875
$ws.Range("L31").Value = 7143.829
$ws.Range("M31").Value = -3201.875
$ws.Range("N31").Value = -7733.829
$ws.Range("H34").Value = 6120.123
$ws.Range("I34").Value = 3496.875
$ws.Range("J34").Value = 7143.829
$ws.Range("K34").Value = 3496.875
$ws.Range("L34").Value = 7143.829
$ws.Range("M34").Value = -3294.875
$ws.Range("N34").Value = -7547.829
$ws.Range("H52").Value = 26181.25
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 26181.25
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 26181.25
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -26769.25
$ws.Range("H58").Value = 16407.354
$ws.Range("I58").Value = 1699.0667
$ws.Range("J58").Value = 28019.158
$ws.Range("K58").Value = 1699.0667
$ws.Range("L58").Value = 28019.158
$ws.Range("M58").Value = -1496.0667
$ws.Range("N58").Value = -28425.158
$ws.Range("H62").Value = 5109.636
$ws.Range("I62").Value = 4550
$ws.Range("J62").Value = 5781.2
$ws.Range("K62").Value = 4550
$ws.Range("L62").Value = 5781.2
$ws.Range("M62").Value = -3926
$ws.Range("N62").Value = -7029.2
$ws.Range("H65").Value = 5109.636
$ws.Range("I65").Value = 4550
$ws.Range("J65").Value = 5781.2
$ws.Range("K65").Value = 22750
$ws.Range("L65").Value = 28906
$ws.Range("M65").Value = -19630
$ws.Range("N65").Value = -35146
$ws.Range("H122").Value = 1774.4736
$ws.Range("I122").Value = 1647.7273
$ws.Range("J122").Value = 1948.75
$ws.Range("K122").Value = 4943.1819
$ws.Range("L122").Value = 5846.25
$ws.Range("M122").Value = -2493.1819
$ws.Range("N122").Value = -10746.25
$ws.Range("H136").Value = 16407.354
$ws.Range("I136").Value = 1699.0667
$ws.Range("J136").Value = 28019.158
$ws.Range("K136").Value = 5097.2001
$ws.Range("L136").Value = 84057.474
$ws.Range("M136").Value = -2547.2001
$ws.Range("N136").Value = -89157.474

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1700.6666
$ws.Range("I36").Value = 549.5
$ws.Range("K36").Value = 1648.5
$ws.Range("M36").Value = -1479.5
$ws.Range("H48").Value = 1050
$ws.Range("I48").Value = 100
$ws.Range("K48").Value = 300
$ws.Range("M48").Value = -50
$ws.Range("H131").Value = 741.39
$ws.Range("J131").Value = 749.6804
$ws.Range("L131").Value = 2249.0412
$ws.Range("N131").Value = -12329.0412
$ws.Range("H132").Value = 769.75
$ws.Range("I132").Value = 769.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6927.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4397.75
$ws.Range("N132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20005600
$ws.Range("J52").Value = 20005600
$ws.Range("L52").Value = 20005600
$ws.Range("N52").Value = -20006118
$ws.Range("H70").Value = 2844977.8
$ws.Range("I70").Value = 4090.818
$ws.Range("J70").Value = 5685864.5
$ws.Range("K70").Value = 4090.818
$ws.Range("L70").Value = 5685864.5
$ws.Range("M70").Value = -3820.818
$ws.Range("N70").Value = -5686404.5
$ws.Range("H73").Value = 2844977.8
$ws.Range("I73").Value = 4090.818
$ws.Range("J73").Value = 5685864.5
$ws.Range("K73").Value = 4090.818
$ws.Range("L73").Value = 5685864.5
$ws.Range("M73").Value = -3154.818
$ws.Range("N73").Value = -5687736.5
$ws.Range("H126").Value = 4814.0713
$ws.Range("I126").Value = 4692.3335
$ws.Range("J126").Value = 5033.2
$ws.Range("K126").Value = 14077.0005
$ws.Range("L126").Value = 15099.6
$ws.Range("M126").Value = -11607.0005
$ws.Range("N126").Value = -20039.6
$ws.Range("H132").Value = 6077012.5
$ws.Range("I132").Value = 9775681
$ws.Range("J132").Value = 66676.5
$ws.Range("K132").Value = 29327043
$ws.Range("L132").Value = 200029.5
$ws.Range("M132").Value = -29324513
$ws.Range("N132").Value = -205089.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 283.33334
$ws.Range("I55").Value = 280
$ws.Range("J55").Value = 295
$ws.Range("K55").Value = 280
$ws.Range("L55").Value = 295
$ws.Range("M55").Value = -107
$ws.Range("N55").Value = -641
$ws.Range("H93").Value = 1785.7142
$ws.Range("I93").Value = 1785.7142
$ws.Range("K93").Value = 1785.7142
$ws.Range("M93").Value = -537.7141999999999
$ws.Range("H122").Value = 1997367.4
$ws.Range("I122").Value = 2218519.2
$ws.Range("K122").Value = 6655557.600000001
$ws.Range("M122").Value = -6653107.600000001
$ws.Range("H132").Value = 3676.923
$ws.Range("J132").Value = 5114.2856
$ws.Range("L132").Value = 15342.8568
$ws.Range("N132").Value = -20402.8568
$ws.Range("H137").Value = 49986.668
$ws.Range("J137").Value = 49986.668
$ws.Range("L137").Value = 49986.668
$ws.Range("N137").Value = -60186.668

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 18000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 18000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 18000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -18584
$ws.Range("H100").Value = 541
$ws.Range("J100").Value = 500
$ws.Range("L100").Value = 1000
$ws.Range("N100").Value = -2082
$ws.Range("H122").Value = 2100
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H136").Value = 52636960
$ws.Range("I136").Value = 111114170
$ws.Range("J136").Value = 7470.5
$ws.Range("K136").Value = 333342510
$ws.Range("L136").Value = 22411.5
$ws.Range("M136").Value = -333339960
$ws.Range("N136").Value = -27511.5

Write-Output "Applied changes: modify=243 add=5 remove=3"
